$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-15 07:56:23"
$wsZhCn.Range("G3").Value = "2016-01-15 07:57:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-15 07:56:33"
$wsDeDe.Range("G3").Value = "2016-01-15 07:57:19"
